$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the data area so the shared-string table is
# rebuilt in the exact order the target file expects (the engine appends newly
# introduced strings in the order cells are written, and compacts/removes any
# strings that end up unreferenced).
$ws.Range("A1:Q27").ClearContents()

# Row 28-31 are brand new rows; clone the bold/bordered label style used by
# A3:H27 onto them before filling values in, so the cells pick up style index 1.
$ws.Range("A27:H27").Copy()
$ws.Range("A28:H31").PasteSpecial(-4122)

# Header row (row 2) - preserve original shared-string order (anchor score.. name)
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("A2").Value = "name"

# Mirror header for the positive block (same strings, order is irrelevant here)
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("J2").Value = "name"

# Negative-word rows (column A) top-to-bottom -> this ordering is exactly what
# the target shared-string table expects for indices 8..36
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9782608695652174
$ws.Range("C3").Value = 45
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 1
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.7954545454545454
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 9
$ws.Range("A5").Value = "however"
$ws.Range("B5").Value = 0.765625
$ws.Range("C5").Value = 49
$ws.Range("D5").Value = 49
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 15
$ws.Range("A6").Value = "poor"
$ws.Range("B6").Value = 0.7323943661971831
$ws.Range("C6").Value = 52
$ws.Range("D6").Value = 52
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 19
$ws.Range("A7").Value = "broke"
$ws.Range("B7").Value = 0.7087378640776699
$ws.Range("C7").Value = 146
$ws.Range("D7").Value = 146
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 60
$ws.Range("A8").Value = "disappointed"
$ws.Range("B8").Value = 0.6989247311827957
$ws.Range("C8").Value = 130
$ws.Range("D8").Value = 130
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 56
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.668918918918919
$ws.Range("C9").Value = 99
$ws.Range("D9").Value = 99
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 49
$ws.Range("A10").Value = "junk"
$ws.Range("B10").Value = 0.6545454545454545
$ws.Range("C10").Value = 36
$ws.Range("D10").Value = 36
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 19
$ws.Range("A11").Value = "smaller"
$ws.Range("B11").Value = 0.5798319327731093
$ws.Range("C11").Value = 69
$ws.Range("D11").Value = 69
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 50
$ws.Range("A12").Value = "small"
$ws.Range("B12").Value = 0.4927536231884058
$ws.Range("C12").Value = 170
$ws.Range("D12").Value = 170
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 175
$ws.Range("A13").Value = "plastic"
$ws.Range("B13").Value = 0.4645669291338583
$ws.Range("C13").Value = 59
$ws.Range("D13").Value = 59
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 68
$ws.Range("A14").Value = "broken"
$ws.Range("B14").Value = 0.4457831325301205
$ws.Range("C14").Value = 37
$ws.Range("D14").Value = 37
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 46
$ws.Range("A15").Value = "apart"
$ws.Range("B15").Value = 0.3894736842105263
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = 37
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 58
$ws.Range("A16").Value = "difficult"
$ws.Range("B16").Value = 0.3483146067415731
$ws.Range("C16").Value = 31
$ws.Range("D16").Value = 31
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 58
$ws.Range("A17").Value = "thought"
$ws.Range("B17").Value = 0.3267326732673267
$ws.Range("C17").Value = 66
$ws.Range("D17").Value = 66
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 136
$ws.Range("A18").Value = "though"
$ws.Range("B18").Value = 0.3162393162393162
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 80
$ws.Range("A19").Value = "ok"
$ws.Range("B19").Value = 0.296875
$ws.Range("C19").Value = 38
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 90
$ws.Range("A20").Value = "cheap"
$ws.Range("B20").Value = 0.2890995260663507
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = 61
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 150
$ws.Range("A21").Value = "size"
$ws.Range("B21").Value = 0.2371134020618557
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 148
$ws.Range("A22").Value = "item"
$ws.Range("B22").Value = 0.1920289855072464
$ws.Range("C22").Value = 53
$ws.Range("D22").Value = 53
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 223
$ws.Range("A23").Value = "would"
$ws.Range("B23").Value = 0.1780415430267062
$ws.Range("C23").Value = 120
$ws.Range("D23").Value = 120
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 554
$ws.Range("A24").Value = "work"
$ws.Range("B24").Value = 0.1677215189873418
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 263
$ws.Range("A25").Value = "money"
$ws.Range("B25").Value = 0.1550632911392405
$ws.Range("C25").Value = 49
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 267
$ws.Range("A26").Value = "better"
$ws.Range("B26").Value = 0.1542056074766355
$ws.Range("C26").Value = 33
$ws.Range("D26").Value = 33
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 181
$ws.Range("A27").Value = "hard"
$ws.Range("B27").Value = 0.15
$ws.Range("C27").Value = 30
$ws.Range("D27").Value = 30
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 170
$ws.Range("A28").Value = "product"
$ws.Range("B28").Value = 0.1233480176211454
$ws.Range("C28").Value = 56
$ws.Range("D28").Value = 56
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 398
$ws.Range("A29").Value = "price"
$ws.Range("B29").Value = 0.1152737752161383
$ws.Range("C29").Value = 40
$ws.Range("D29").Value = 41
$ws.Range("E29").Value = 0.02
$ws.Range("F29").Value = 0.98
$ws.Range("G29").Value = $true
$ws.Range("H29").Value = 307
$ws.Range("A30").Value = "use"
$ws.Range("B30").Value = 0.08767123287671233
$ws.Range("C30").Value = 32
$ws.Range("D30").Value = 32
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 333
$ws.Range("A31").Value = "like"
$ws.Range("B31").Value = 0.06578947368421052
$ws.Range("C31").Value = 40
$ws.Range("D31").Value = 40
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = $false
$ws.Range("H31").Value = 568

# Top label for the negative block, introduced right after the word list
$ws.Range("A1").Value = "negative"

# Positive-word rows (column J) top-to-bottom -> indices 38..49
$ws.Range("J3").Value = "wonderful"
$ws.Range("K3").Value = 0.7857142857142857
$ws.Range("L3").Value = 44
$ws.Range("M3").Value = 44
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 12
$ws.Range("J4").Value = "awesome"
$ws.Range("K4").Value = 0.7538461538461538
$ws.Range("L4").Value = 49
$ws.Range("M4").Value = 49
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 16
$ws.Range("J5").Value = "favorite"
$ws.Range("K5").Value = 0.6666666666666666
$ws.Range("L5").Value = 62
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 31
$ws.Range("J6").Value = "classic"
$ws.Range("K6").Value = 0.5849056603773585
$ws.Range("L6").Value = 31
$ws.Range("M6").Value = 31
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 22
$ws.Range("J7").Value = "excellent"
$ws.Range("K7").Value = 0.46875
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 34
$ws.Range("J8").Value = "great"
$ws.Range("K8").Value = 0.3418032786885246
$ws.Range("L8").Value = 417
$ws.Range("M8").Value = 417
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 803
$ws.Range("J9").Value = "love"
$ws.Range("K9").Value = 0.302725968436155
$ws.Range("L9").Value = 211
$ws.Range("M9").Value = 211
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 486
$ws.Range("J10").Value = "loves"
$ws.Range("K10").Value = 0.2448132780082987
$ws.Range("L10").Value = 118
$ws.Range("M10").Value = 118
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 364
$ws.Range("J11").Value = "perfect"
$ws.Range("K11").Value = 0.1927710843373494
$ws.Range("L11").Value = 32
$ws.Range("M11").Value = 32
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 134
$ws.Range("J12").Value = "loved"
$ws.Range("K12").Value = 0.1620795107033639
$ws.Range("L12").Value = 53
$ws.Range("M12").Value = 53
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 274
$ws.Range("J13").Value = "fun"
$ws.Range("K13").Value = 0.08413672217353199
$ws.Range("L13").Value = 96
$ws.Range("M13").Value = 96
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 1045
$ws.Range("J14").Value = "game"
$ws.Range("K14").Value = 0.03439325113562622
$ws.Range("L14").Value = 53
$ws.Range("M14").Value = 53
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 1488

# Top label for the positive block (last new string -> index 50)
$ws.Range("J1").Value = "positive"

